# Update the repayment figures for the Hansyah_S2l / S2 collector sheet.
# (Reflects a refreshed daily repayment export: updated Talk_time totals for
# every collector, plus a handful of collectors who picked up new
# repayments since the previous upload.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: "(3)" -> "(4)" to reflect the new upload.
$ws.Name = "repayment_20250915_20250915 (4)"

# Helper: write a value as plain text (matches the source data, where the
# amount/rate columns are stored as text strings like "885,787.00" or
# "0.47" rather than numbers) without leaving a lingering number-format
# override on the cell.
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# --- Talk_time (column H) refreshed for every collector row ---
$ws.Range("H2").Value = 348
$ws.Range("H3").Value = 217
$ws.Range("H4").Value = 537
$ws.Range("H5").Value = 617
$ws.Range("H6").Value = 597
$ws.Range("H7").Value = 443
$ws.Range("H8").Value = 922
$ws.Range("H9").Value = 469
$ws.Range("H10").Value = 339
$ws.Range("H11").Value = 367
$ws.Range("H12").Value = 596
$ws.Range("H13").Value = 931
$ws.Range("H14").Value = 35
$ws.Range("H15").Value = 361
$ws.Range("H16").Value = 381
$ws.Range("H17").Value = 514
$ws.Range("H18").Value = 1.4

# --- Row 8: Annisa Putri Restu picked up an extra repayment ---
$ws.Range("D8").Value = 2
Set-TextValue $ws.Range("E8") "885,787.00"
Set-TextValue $ws.Range("G8") "0.47"
$ws.Range("J8").Value = 1
Set-TextValue $ws.Range("L8") "3.33"

# --- Row 10: Azizah Rahmawati now has 2 repayment collections ---
$ws.Range("D10").Value = 2
Set-TextValue $ws.Range("E10") "311,014.00"
Set-TextValue $ws.Range("G10") "0.18"

# --- Row 13: Romli now has 1 repayment collection ---
$ws.Range("D13").Value = 1
Set-TextValue $ws.Range("E13") "50,000.00"
Set-TextValue $ws.Range("G13") "0.03"

# --- Row 16: Yandi Nugraha now has 1 repayment collection ---
$ws.Range("D16").Value = 1
Set-TextValue $ws.Range("E16") "982,315.00"
Set-TextValue $ws.Range("G16") "0.76"
